$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check" header text in F1
$ws.Range("F1").Value = "Last status check on: 24.02.2022 10:45"

# D6: convert from text "+0.39" to a real number 0.39
$ws.Range("D6").Value = 0.39

# E6: convert from text date string to a real date/time serial,
# matching the date formatting used by the other rows in column E
$ws.Range("E6").Value = 44616.43773148148
$ws.Range("E6").NumberFormat = "YYYY-MM-DD HH:MM:SS"
